$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Costs and Revenues")
$ws.Range("B2").Value = 73262.59249872465
$ws.Range("C2").Value = 0
$ws.Range("D2").Value = 9095.293897561472
$ws.Range("E2").Value = 2160
$ws.Range("F2").Value = 21838.00813806339

$ws = $wb.Worksheets.Item("Capacities")
$ws.Range("B3").Value = 0
$ws.Range("C3").Value = 94
$ws.Range("B4").Value = 0
$ws.Range("C4").Value = 150
$ws.Range("D4").Value = 11

$ws = $wb.Worksheets.Item("PV Dispatch")
$ws.Range("G2").Value = 18.8
$ws.Range("H2").Value = 37.6
$ws.Range("I2").Value = 47
$ws.Range("J2").Value = 56.4
$ws.Range("K2").Value = 65.8
$ws.Range("L2").Value = 75.2
$ws.Range("M2").Value = 84.59999999999999
$ws.Range("N2").Value = 94
$ws.Range("O2").Value = 84.59999999999999
$ws.Range("P2").Value = 75.2
$ws.Range("Q2").Value = 65.8
$ws.Range("R2").Value = 47
$ws.Range("S2").Value = 28.2
$ws.Range("T2").Value = 18.8
$ws.Range("I3").Value = 0
$ws.Range("J3").Value = 56.4
$ws.Range("K3").Value = 75.2
$ws.Range("L3").Value = 25.89372513008871
$ws.Range("M3").Value = 94
$ws.Range("N3").Value = 26
$ws.Range("O3").Value = 65.8
$ws.Range("P3").Value = 47
$ws.Range("Q3").Value = 47
$ws.Range("R3").Value = 28.2
$ws.Range("S3").Value = 18.8
$ws.Range("J4").Value = 9.4
$ws.Range("L4").Value = 65.8
$ws.Range("M4").Value = 75.2
$ws.Range("N4").Value = 71.38312417100299
$ws.Range("O4").Value = 0
$ws.Range("P4").Value = 37.6
$ws.Range("R4").Value = 9.4

$ws = $wb.Worksheets.Item("Battery Input")
$ws.Range("G2").Value = 11
$ws.Range("H2").Value = 24.6
$ws.Range("I2").Value = 15.8
$ws.Range("J2").Value = 17.4
$ws.Range("K2").Value = 39.8
$ws.Range("L2").Value = 54.4
$ws.Range("M2").Value = 61.2
$ws.Range("N2").Value = 68
$ws.Range("O2").Value = 92.33427201305987
$ws.Range("P2").Value = 46.6
$ws.Range("Q2").Value = 39.8
$ws.Range("R2").Value = 13.2
$ws.Range("I3").Value = 0
$ws.Range("J3").Value = 56.4
$ws.Range("K3").Value = 75.2
$ws.Range("L3").Value = 25.89372513008871
$ws.Range("M3").Value = 70.59999999999999
$ws.Range("N3").Value = 0
$ws.Range("O3").Value = 65.8
$ws.Range("P3").Value = 18.4
$ws.Range("Q3").Value = 21
$ws.Range("R3").Value = 28.2
$ws.Range("J4").Value = 9.4
$ws.Range("L4").Value = 65.8
$ws.Range("M4").Value = 51.8
$ws.Range("N4").Value = 71.38312417100299
$ws.Range("O4").Value = 0
$ws.Range("P4").Value = 37.6
$ws.Range("R4").Value = 9.4

$ws = $wb.Worksheets.Item("Battery Output")
$ws.Range("S2").Value = 0
$ws.Range("T2").Value = 0
$ws.Range("S3").Value = 22.8

$ws = $wb.Worksheets.Item("State of Charge")
$ws.Range("B2").Value = 179.0909090909091
$ws.Range("C2").Value = 159.3939393939394
$ws.Range("D2").Value = 146.2626262626262
$ws.Range("E2").Value = 133.1313131313131
$ws.Range("F2").Value = 120
$ws.Range("G2").Value = 130.89
$ws.Range("H2").Value = 155.244
$ws.Range("I2").Value = 170.886
$ws.Range("J2").Value = 188.112
$ws.Range("K2").Value = 227.514
$ws.Range("L2").Value = 281.37
$ws.Range("M2").Value = 341.958
$ws.Range("N2").Value = 409.278
$ws.Range("O2").Value = 500.6889292929293
$ws.Range("P2").Value = 546.8229292929293
$ws.Range("Q2").Value = 586.2249292929293
$ws.Range("R2").Value = 599.2929292929293
$ws.Range("S2").Value = 599.2929292929293
$ws.Range("T2").Value = 599.2929292929293
$ws.Range("U2").Value = 481.1111111111111
$ws.Range("V2").Value = 382.6262626262626
$ws.Range("W2").Value = 303.8383838383838
$ws.Range("X2").Value = 251.3131313131313
$ws.Range("Y2").Value = 211.9191919191919
$ws.Range("B3").Value = 172.5252525252525
$ws.Range("C3").Value = 152.8282828282828
$ws.Range("D3").Value = 139.6969696969697
$ws.Range("E3").Value = 139.6969696969697
$ws.Range("F3").Value = 139.6969696969697
$ws.Range("G3").Value = 120
$ws.Range("H3").Value = 120
$ws.Range("I3").Value = 120
$ws.Range("J3").Value = 175.836
$ws.Range("K3").Value = 250.284
$ws.Range("L3").Value = 275.9187878787878
$ws.Range("M3").Value = 345.8127878787878
$ws.Range("N3").Value = 345.8127878787878
$ws.Range("O3").Value = 410.9547878787878
$ws.Range("P3").Value = 429.1707878787878
$ws.Range("Q3").Value = 449.9607878787879
$ws.Range("R3").Value = 477.8787878787879
$ws.Range("S3").Value = 454.8484848484849
$ws.Range("T3").Value = 323.5353535353535
$ws.Range("U3").Value = 323.5353535353535
$ws.Range("V3").Value = 323.5353535353535
$ws.Range("W3").Value = 244.7474747474747
$ws.Range("X3").Value = 244.7474747474747
$ws.Range("Y3").Value = 205.3535353535353
$ws.Range("B4").Value = 159.3939393939394
$ws.Range("C4").Value = 139.6969696969697
$ws.Range("D4").Value = 139.6969696969697
$ws.Range("E4").Value = 139.6969696969697
$ws.Range("F4").Value = 139.6969696969697
$ws.Range("G4").Value = 120
$ws.Range("H4").Value = 120
$ws.Range("I4").Value = 120
$ws.Range("J4").Value = 129.306
$ws.Range("K4").Value = 129.306
$ws.Range("L4").Value = 194.448
$ws.Range("M4").Value = 245.73
$ws.Range("N4").Value = 316.399292929293
$ws.Range("O4").Value = 316.399292929293
$ws.Range("P4").Value = 353.623292929293
$ws.Range("Q4").Value = 353.623292929293
$ws.Range("R4").Value = 362.929292929293
$ws.Range("S4").Value = 362.929292929293
$ws.Range("T4").Value = 231.6161616161616
$ws.Range("U4").Value = 231.6161616161616
$ws.Range("V4").Value = 231.6161616161616
$ws.Range("W4").Value = 231.6161616161616
$ws.Range("X4").Value = 231.6161616161616
$ws.Range("Y4").Value = 192.2222222222222

$ws = $wb.Worksheets.Item("Feed in from Type 2")
$ws.Range("O2").Value = 38.93427201305987
$ws.Range("S2").Value = 13.4
$ws.Range("T2").Value = 15.2

$ws = $wb.Worksheets.Item("Feed in from Type 3")
$ws.Range("T2").Value = 18

